# Auto-generated Excel COM-interop script to apply the Lamia_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2391.6
$ws.Range("I18").Value = 2391.6
$ws.Range("K18").Value = 2391.6
$ws.Range("M18").Value = -2107.6
$ws.Range("H32").Value = 8039.222
$ws.Range("I32").Value = 7972
$ws.Range("J32").Value = 8093
$ws.Range("K32").Value = 7972
$ws.Range("L32").Value = 8093
$ws.Range("M32").Value = -7646
$ws.Range("N32").Value = -8745
$ws.Range("H69").Value = 7677.5
$ws.Range("J69").Value = 7771.1177
$ws.Range("L69").Value = 23313.3531
$ws.Range("N69").Value = -25061.3531
$ws.Range("H72").Value = 7677.5
$ws.Range("J72").Value = 7771.1177
$ws.Range("L72").Value = 69940.05929999999
$ws.Range("N72").Value = -78676.05929999999
$ws.Range("H100").Value = 4440.222
$ws.Range("I100").Value = 2083.4
$ws.Range("J100").Value = 7386.25
$ws.Range("K100").Value = 2083.4
$ws.Range("L100").Value = 7386.25
$ws.Range("M100").Value = -1542.4
$ws.Range("N100").Value = -8468.25
$ws.Range("H137").Value = 3087.2922
$ws.Range("I137").Value = 2436.2856
$ws.Range("J137").Value = 3266
$ws.Range("K137").Value = 7308.8568
$ws.Range("L137").Value = 9798
$ws.Range("M137").Value = -4758.8568
$ws.Range("N137").Value = -14898
$ws.Range("H138").Value = 2824.7468
$ws.Range("I138").Value = 1542.625
$ws.Range("J138").Value = 3384.2183
$ws.Range("K138").Value = 4627.875
$ws.Range("L138").Value = 10152.6549
$ws.Range("M138").Value = 512.125
$ws.Range("N138").Value = -20432.6549

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5391.364
$ws.Range("I63").Value = 3272.923
$ws.Range("K63").Value = 3272.923
$ws.Range("M63").Value = -2586.923
$ws.Range("H66").Value = 5391.364
$ws.Range("I66").Value = 3272.923
$ws.Range("K66").Value = 16364.615
$ws.Range("M66").Value = -12932.615
$ws.Range("H110").Value = 4844.88
$ws.Range("I110").Value = 4456.5557
$ws.Range("J110").Value = 5843.4287
$ws.Range("K110").Value = 4456.5557
$ws.Range("L110").Value = 5843.4287
$ws.Range("M110").Value = -2411.5557
$ws.Range("N110").Value = -9933.4287
$ws.Range("H132").Value = 2740.4443
$ws.Range("I132").Value = 1999.12
$ws.Range("J132").Value = 12007
$ws.Range("K132").Value = 5997.36
$ws.Range("L132").Value = 36021
$ws.Range("M132").Value = -3467.36
$ws.Range("N132").Value = -41081

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 17484.5
$ws.Range("I128").Value = 17484.5
$ws.Range("K128").Value = 52453.5
$ws.Range("M128").Value = -49963.5
$ws.Range("H134").Value = 3130.4375
$ws.Range("I134").Value = 1648.0714
$ws.Range("K134").Value = 4944.2142
$ws.Range("M134").Value = -2409.2142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30155.21
$ws.Range("I31").Value = 1643
$ws.Range("K31").Value = 1643
$ws.Range("M31").Value = -1348
$ws.Range("H32").Value = 2105
$ws.Range("I32").Value = 1806.6666
$ws.Range("K32").Value = 1806.6666
$ws.Range("M32").Value = -1490.6666
$ws.Range("H34").Value = 30155.21
$ws.Range("I34").Value = 1643
$ws.Range("K34").Value = 1643
$ws.Range("M34").Value = -1441
$ws.Range("H37").Value = 32000
$ws.Range("J37").Value = 32000
$ws.Range("L37").Value = 32000
$ws.Range("N37").Value = -32214
$ws.Range("H94").Value = 2554.2666
$ws.Range("J94").Value = 3292.875
$ws.Range("L94").Value = 3292.875
$ws.Range("N94").Value = -4194.875
$ws.Range("H134").Value = 2347.074
$ws.Range("I134").Value = 1232.9565
$ws.Range("K134").Value = 3698.8695
$ws.Range("M134").Value = -1163.8695

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2399
$ws.Range("I63").Value = 2399
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 7197
$ws.Range("L63").ClearContents()
$ws.Range("N63").Value = 0
$ws.Range("M63").Value = -6448
$ws.Range("H64").Value = 111118810
$ws.Range("I64").Value = 1000000000
$ws.Range("J64").Value = 8655.25
$ws.Range("K64").Value = 3000000000
$ws.Range("L64").Value = 25965.75
$ws.Range("M64").Value = -2999999730
$ws.Range("N64").Value = -26505.75
$ws.Range("H66").Value = 2399
$ws.Range("I66").Value = 2399
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 21591
$ws.Range("L66").ClearContents()
$ws.Range("N66").Value = 0
$ws.Range("M66").Value = -17847
$ws.Range("H67").Value = 111118810
$ws.Range("I67").Value = 1000000000
$ws.Range("J67").Value = 8655.25
$ws.Range("K67").Value = 3000000000
$ws.Range("L67").Value = 25965.75
$ws.Range("M67").Value = -2999999064
$ws.Range("N67").Value = -27837.75
$ws.Range("H114").Value = 1553.375
$ws.Range("I114").Value = 1131.75
$ws.Range("J114").Value = 1975
$ws.Range("K114").Value = 3395.25
$ws.Range("L114").Value = 5925
$ws.Range("M114").Value = -141.25
$ws.Range("N114").Value = -12433
$ws.Range("H117").Value = 1477.6364
$ws.Range("J117").Value = 3445
$ws.Range("L117").Value = 10335
$ws.Range("N117").Value = -17219
$ws.Range("H132").Value = 4946.067
$ws.Range("I132").Value = 3649
$ws.Range("J132").Value = 6428.4287
$ws.Range("K132").Value = 32841
$ws.Range("L132").Value = 57855.85830000001
$ws.Range("M132").Value = -30311
$ws.Range("N132").Value = -62915.85830000001
$ws.Range("H139").Value = 4747.2144
$ws.Range("I139").Value = 2971.4285
$ws.Range("J139").Value = 6523
$ws.Range("K139").Value = 8914.2855
$ws.Range("L139").Value = 19569
$ws.Range("M139").Value = -3774.2855
$ws.Range("N139").Value = -29849
$ws.Range("H141").Value = 6675.9546
$ws.Range("I141").Value = 1404.0834
$ws.Range("J141").Value = 13002.2
$ws.Range("K141").Value = 4212.2502
$ws.Range("L141").Value = 39006.60000000001
$ws.Range("M141").Value = 967.7497999999996
$ws.Range("N141").Value = -49366.60000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 42499.5
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 42499.5
$ws.Range("K44").Value = 0
$ws.Range("L44").ClearContents()
$ws.Range("M44").Value = 42499.5
$ws.Range("N44").Value = -43691.5
$ws.Range("H52").Value = 25010000
$ws.Range("I52").Value = 50000000
$ws.Range("K52").Value = 50000000
$ws.Range("M52").Value = -49999741
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").ClearContents()
$ws.Range("N101").Value = 0
$ws.Range("H122").Value = 8678.5
$ws.Range("I122").Value = 10077.105
$ws.Range("J122").Value = 7115.353
$ws.Range("K122").Value = 30231.315
$ws.Range("L122").Value = 21346.059
$ws.Range("M122").Value = -27781.315
$ws.Range("N122").Value = -26246.059
$ws.Range("H126").Value = 5346.727
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 2549.6667
$ws.Range("I132").Value = 1340.2858
$ws.Range("J132").Value = 4242.8
$ws.Range("K132").Value = 4020.8574
$ws.Range("L132").Value = 12728.4
$ws.Range("M132").Value = -1490.8574
$ws.Range("N132").Value = -17788.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6194.7095
$ws.Range("I7").Value = 4304.56
$ws.Range("J7").Value = 14070.333
$ws.Range("K7").Value = 4304.56
$ws.Range("L7").Value = 14070.333
$ws.Range("M7").Value = -4192.56
$ws.Range("N7").Value = -14294.333
$ws.Range("H22").Value = 3013.6924
$ws.Range("I22").Value = 1780.4375
$ws.Range("J22").Value = 4986.9
$ws.Range("K22").Value = 1780.4375
$ws.Range("L22").Value = 4986.9
$ws.Range("M22").Value = -1485.4375
$ws.Range("N22").Value = -5576.9
$ws.Range("H27").Value = 3013.6924
$ws.Range("I27").Value = 1780.4375
$ws.Range("J27").Value = 4986.9
$ws.Range("K27").Value = 1780.4375
$ws.Range("L27").Value = 4986.9
$ws.Range("M27").Value = -1673.4375
$ws.Range("N27").Value = -5200.9
$ws.Range("H32").Value = 3138.125
$ws.Range("I32").Value = 3086.4285
$ws.Range("K32").Value = 3086.4285
$ws.Range("M32").Value = -2769.4285
$ws.Range("H82").Value = 4888.9287
$ws.Range("I82").Value = 1224
$ws.Range("K82").Value = 1224
$ws.Range("M82").Value = -863
$ws.Range("H85").Value = 4888.9287
$ws.Range("I85").Value = 1224
$ws.Range("K85").Value = 1224
$ws.Range("M85").Value = 24
$ws.Range("H100").Value = 12637.1
$ws.Range("I100").Value = 8284.5
$ws.Range("K100").Value = 8284.5
$ws.Range("M100").Value = -7743.5
$ws.Range("H126").Value = 6194.7095
$ws.Range("I126").Value = 4304.56
$ws.Range("J126").Value = 14070.333
$ws.Range("K126").Value = 12913.68
$ws.Range("L126").Value = 42210.999
$ws.Range("M126").Value = -10443.68
$ws.Range("N126").Value = -47150.999
$ws.Range("H132").Value = 6092.8184
$ws.Range("I132").Value = 4476.077
$ws.Range("J132").Value = 8428.111000000001
$ws.Range("K132").Value = 13428.231
$ws.Range("L132").Value = 25284.333
$ws.Range("M132").Value = -10898.231
$ws.Range("N132").Value = -30344.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").ClearContents()
$ws.Range("N60").Value = 0
$ws.Range("H132").Value = 2668.1428
$ws.Range("I132").Value = 1576.9375
$ws.Range("K132").Value = 4730.8125
$ws.Range("M132").Value = -2200.8125
$ws.Range("H136").Value = 3083.7693
$ws.Range("I136").Value = 2047.5238
$ws.Range("K136").Value = 6142.5714
$ws.Range("M136").Value = -3592.5714
